$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-16 (columns A: Name, B: Position, C: Team)
$data = @(
    @(2,  "Coby White",          "PG,SG",    "Chicago Bulls"),
    @(3,  "Devin Booker",        "PG,SG",    "Phoenix Suns"),
    @(5,  "Trae Young",          "PG",       "Atlanta Hawks"),
    @(8,  "Shaedon Sharpe",      "SG,SF",    "Portland Trail Blazers"),
    @(9,  "Goga Bitadze",        "C",        "Orlando Magic"),
    @(10, "Alperen Sengün",      "C",        "Houston Rockets"),
    @(11, "LeBron James",        "SF,PF",    "Los Angeles Lakers"),
    @(13, "Stephon Castle",      "PG,SG",    "San Antonio Spurs"),
    @(14, "Devin Vassell",       "SG,SF",    "San Antonio Spurs"),
    @(15, "Walker Kessler",      "C",        "Utah Jazz"),
    @(16, "Norman Powell",       "SG,SF",    "LA Clippers")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
